$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 501-502; existing rows 501-574 shift down to 503-576
$ws.Rows("501:502").Insert()

# --- New row 501 ---
$ws.Range("A501").Value = 7
$ws.Range("B501").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C501").Value = "Ñuble"
$ws.Range("D501").Value = [DateTime]::FromOADate(45180)
$ws.Range("E501").Value = 16
$ws.Range("F501").Value = 100112009
$ws.Range("G501").Value = "Acelga"
$ws.Range("H501").Value = "Sin especificar"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 300
$ws.Range("K501").Value = 600
$ws.Range("L501").Value = 700
$ws.Range("M501").Value = 650
$ws.Range("N501").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O501").Value = "Región de Ñuble"
$ws.Range("P501").Value = 650
$ws.Range("Q501").Value = 1
$ws.Range("R501").Value = "Hortaliza"

# --- New row 502 ---
$ws.Range("A502").Value = 7
$ws.Range("B502").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C502").Value = "Ñuble"
$ws.Range("D502").Value = [DateTime]::FromOADate(45180)
$ws.Range("E502").Value = 16
$ws.Range("F502").Value = 100112009
$ws.Range("G502").Value = "Acelga"
$ws.Range("H502").Value = "Sin especificar"
$ws.Range("I502").Value = "Segunda"
$ws.Range("J502").Value = 250
$ws.Range("K502").Value = 500
$ws.Range("L502").Value = 500
$ws.Range("M502").Value = 500
$ws.Range("N502").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O502").Value = "Región de Ñuble"
$ws.Range("P502").Value = 500
$ws.Range("Q502").Value = 1
$ws.Range("R502").Value = "Hortaliza"

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
